$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Row 9: device changed from MZX125 / 0.080 to MZX250 / 0.105
$ws.Range("F9").Value = "'0.105"
$ws.Range("A9").Value = "MZX250"

# B4: new device location/type text (was empty)
$ws.Range("B4").Value = "NGC-601/T1465 OR TC-217"

# Update the active selection to match the authored state
$ws.Range("B4").Select() | Out-Null
